# Cronograma.xlsx - "Se actualizaron fechas de trabajo y días para próximas evaluaciones"
#
# The sheet is a weekly schedule (cols A/I = time-of-day, B..G and J..O = days).
# Columns N/O (and occasionally M) hold the "Descanso" vs "Espacio de repaso"
# marker for a given slot, while column M sometimes holds "Trabajo" instead.
# This edit shifts several slots from "Trabajo" -> "Descanso" and several
# "Descanso" -> "Espacio de repaso" (making room for upcoming evaluations),
# and moves the active selection/scroll position in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cells whose value+format we reuse (via copy/paste-special) so the
# edited cells pick up the exact same existing style record (fill/border)
# instead of Excel fabricating a brand new one for each cell.
$srcEspacioDeRepaso = $ws.Range("C18")   # "Espacio de repaso" (red fill) style
$srcDescanso        = $ws.Range("O13")   # "Descanso" (bordered, no fill) style

# N13: Descanso -> Espacio de repaso
$ws.Range("N13").Value = "Espacio de repaso"
$srcEspacioDeRepaso.Copy()
$ws.Range("N13").PasteSpecial(-4122)

# M14: Trabajo -> Descanso
$ws.Range("M14").Value = "Descanso"
$srcDescanso.Copy()
$ws.Range("M14").PasteSpecial(-4122)

# N14: Descanso -> Espacio de repaso
$ws.Range("N14").Value = "Espacio de repaso"
$srcEspacioDeRepaso.Copy()
$ws.Range("N14").PasteSpecial(-4122)

# N15: Descanso -> Espacio de repaso
$ws.Range("N15").Value = "Espacio de repaso"
$srcEspacioDeRepaso.Copy()
$ws.Range("N15").PasteSpecial(-4122)

# M16: Trabajo -> Descanso
$ws.Range("M16").Value = "Descanso"
$srcDescanso.Copy()
$ws.Range("M16").PasteSpecial(-4122)

# O16: Descanso -> Espacio de repaso
$ws.Range("O16").Value = "Espacio de repaso"
$srcEspacioDeRepaso.Copy()
$ws.Range("O16").PasteSpecial(-4122)

# M18: Trabajo -> Descanso
$ws.Range("M18").Value = "Descanso"
$srcDescanso.Copy()
$ws.Range("M18").PasteSpecial(-4122)

# Update the window scroll position (top-left visible cell E2) and selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 5
$ws.Range("M14").Select()
